$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure these cells stay stored as text so exact formatting (trailing zeros, % sign) is preserved
$cells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "E17", "E18", "D19", "E19", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($c in $cells) { $ws.Range($c).NumberFormat = "@" }

$ws.Range("D2").Value = "305.45"
$ws.Range("E2").Value = "0.40%"
$ws.Range("D3").Value = "35.81"
$ws.Range("E3").Value = "-0.92%"
$ws.Range("D4").Value = "4.968"
$ws.Range("E4").Value = "-1.77%"
$ws.Range("D5").Value = "0.08076"
$ws.Range("E5").Value = "-0.65%"
$ws.Range("D6").Value = "1.899"
$ws.Range("E6").Value = "-3.38%"
$ws.Range("D7").Value = "4.148"
$ws.Range("E7").Value = "1.95%"
$ws.Range("D8").Value = "7.879"
$ws.Range("E8").Value = "0.69%"
$ws.Range("D9").Value = "0.9297"
$ws.Range("E10").Value = "-20.30%"
$ws.Range("D11").Value = "0.1905"
$ws.Range("E11").Value = "-0.18%"
$ws.Range("D12").Value = "0.09222"
$ws.Range("E12").Value = "2.04%"
$ws.Range("D13").Value = "0.03506"
$ws.Range("E13").Value = "1.59%"
$ws.Range("D14").Value = "0.09911"
$ws.Range("E14").Value = "0.71%"
$ws.Range("D15").Value = "0.001425"
$ws.Range("E15").Value = "0.29%"
$ws.Range("D16").Value = "0.006336"
$ws.Range("E16").Value = "7.88%"
$ws.Range("E17").Value = "1.86%"
$ws.Range("E18").Value = "4.18%"
$ws.Range("D19").Value = "0.3445"
$ws.Range("E19").Value = "-0.18%"
$ws.Range("E21").Value = "2.96%"
$ws.Range("D22").Value = "0.2531"
$ws.Range("E22").Value = "5.93%"
$ws.Range("D23").Value = "0.04415"
$ws.Range("E23").Value = "-1.15%"
$ws.Range("D24").Value = "0.001235"
$ws.Range("E24").Value = "2.63%"
$ws.Range("D25").Value = "0.004714"
$ws.Range("E25").Value = "-2.62%"
$ws.Range("D26").Value = "0.0001300"
$ws.Range("E26").Value = "6.22%"
$ws.Range("D27").Value = "0.0003129"
$ws.Range("E27").Value = "3.79%"
$ws.Range("D39").Value = "0.01951"
$ws.Range("E39").Value = "0.87%"
$ws.Range("D40").Value = "0.05200"
$ws.Range("E40").Value = "8.19%"
$ws.Range("D41").Value = "0.007547"
$ws.Range("E41").Value = "3.33%"
$ws.Range("D42").Value = "0.01016"
$ws.Range("E42").Value = "-4.32%"
$ws.Range("D43").Value = "0.1370"
$ws.Range("E43").Value = "1.26%"
$ws.Range("D44").Value = "0.002101"
$ws.Range("E44").Value = "0.02%"
$ws.Range("E45").Value = "-0.37%"
$ws.Range("D46").Value = "0.00006345"
$ws.Range("E46").Value = "4.31%"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").Value = "0.47%"
$ws.Range("D48").Value = "63.57"
$ws.Range("E48").Value = "-1.70%"
$ws.Range("D49").Value = "0.001659"
$ws.Range("E49").Value = "-0.02%"
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").Value = "0.47%"
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").Value = "0.47%"
